$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.186.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.23%  "

$ws.Range("D3").Value = "'2.665.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'596.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").Value = "'163.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.544"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'2.665.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.80%  "

$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").Value = "'0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").Value = "'5.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("D14").Value = "'27.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").Value = "'3.164.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").Value = "'0.0000183"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.63%  "

$ws.Range("D17").Value = "'67.209.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("D18").Value = "'2.665.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.47%  "

$ws.Range("D19").Value = "'11.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.15%  "

$ws.Range("D20").Value = "'361.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.65%  "

$ws.Range("D21").Value = "'7.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.82%  "

$ws.Range("D22").Value = "'4.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.15%  "

$ws.Range("D23").Value = "'4.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.01%  "

$ws.Range("D24").Value = "'2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.37%  "

$ws.Range("D25").Value = "'71.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.18%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "'10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "'2.817.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -2.88%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").Value = "'549.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.06%  "

$ws.Range("D32").Value = "'7.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.52%  "

$ws.Range("E33").Value = "  -3.72%  "

$ws.Range("D34").Value = "'1.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.47%  "

$ws.Range("D35").Value = "'0.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.67%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "'1.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.32%  "

$ws.Range("D38").Value = "'19.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.89%  "

$ws.Range("D39").Value = "'156.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.20%  "

$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("D41").Value = "'1.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.57%  "

$ws.Range("D42").Value = "'5.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.52%  "

$ws.Range("D43").Value = "'17.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("E44").Value = "  -5.09%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "'40.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").Value = "'0.0₆0298"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.68%  "

$ws.Range("D48").Value = "'0.583"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.81%  "

$ws.Range("D49").Value = "'152.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.98%  "

$ws.Range("D50").Value = "'3.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.18%  "

$ws.Range("D51").Value = "'1.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.04%  "
